# Data update 2023-05-31: append the FTSE 100 ticker list (rows A2:A61)
# again at the bottom of the sheet (rows A2367:A2426).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @(
    "AAF", "ABDN", "ABF", "ANTO", "AUTO", "AV", "BARC", "BATS", "BDEV", "BEZ",
    "BF.B", "BKG", "BNZL", "BRBY", "BRK.B", "BT-A", "CCH", "CRDA", "DCC", "DGE",
    "ENT", "EXPN", "FCIT", "FRAS", "GLEN", "HLMA", "HSBA", "HSX", "IMB", "INF",
    "ITRK", "JMAT", "KGF", "LGEN", "LLOY", "LSEG", "MNDI", "MNG", "OCDO", "PHNX",
    "PSON", "REL", "RMV", "RR", "RS1", "SBRY", "SDR", "SGRO", "SKG", "SMDS",
    "SMT", "SN", "SPX", "SSE", "STAN", "STJ", "ULVR", "UU", "WEIR", "WTB"
)

$startRow = 2367
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 1).Value = $tickers[$i]
}
